$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report timestamp in title (A1)
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 18:50"

# Per-row updates (country reassignments from shared-string reorder + updated case counts)
$ws.Range("B4").Value = 228727
$ws.Range("C4").Value = 13724
$ws.Range("D4").Value = 10280
$ws.Range("E4").Value = 213074
$ws.Range("F4").Value = 5421
$ws.Range("G4").Value = 271
$ws.Range("H4").Value = 5373

$ws.Range("B13").Value = 18135
$ws.Range("C13").Value = 2456
$ws.Range("D13").Value = 415
$ws.Range("E13").Value = 17364
$ws.Range("F13").Value = 783
$ws.Range("G13").Value = 79
$ws.Range("H13").Value = 356

$ws.Range("E32").Value = 2364
$ws.Range("G32").Value = 15
$ws.Range("H32").Value = 107

$ws.Range("A42").Value = "Grecia"
$ws.Range("B42").Value = 1544
$ws.Range("C42").Value = 129
$ws.Range("D42").Value = 61
$ws.Range("E42").Value = 1430
$ws.Range("F42").Value = 91
$ws.Range("H42").Value = 53

$ws.Range("A43").Value = "Finlandia"
$ws.Range("B43").Value = 1518
$ws.Range("C43").Value = 72
$ws.Range("D43").Value = 300
$ws.Range("E43").Value = 1199
$ws.Range("F43").Value = 62
$ws.Range("H43").Value = 19

$ws.Range("D73").Value = 46
$ws.Range("E73").Value = 432

$ws.Range("E108").Value = 126
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 4

$ws.Range("A129").Value = "Monaco"
$ws.Range("B129").Value = 60
$ws.Range("C129").Value = 5
$ws.Range("D129").Value = 2
$ws.Range("E129").Value = 57
$ws.Range("F129").Value = 2
$ws.Range("H129").Value = 1

$ws.Range("A130").Value = "Madagascar"
$ws.Range("B130").Value = 59
$ws.Range("D130").Value = 0
$ws.Range("E130").Value = 59
$ws.Range("F130").Value = 6
$ws.Range("H130").Value = 0

$ws.Range("A131").Value = "Banglades"
$ws.Range("B131").Value = 56
$ws.Range("C131").Value = 2
$ws.Range("D131").Value = 25
$ws.Range("E131").Value = 25
$ws.Range("F131").Value = 1
$ws.Range("H131").Value = 6

$ws.Range("B134").Value = 47
$ws.Range("C134").Value = 8
$ws.Range("E134").Value = 34

$ws.Range("A144").Value = "Mali"
$ws.Range("C144").Value = 5
$ws.Range("D144").Value = 0
$ws.Range("E144").Value = 33
$ws.Range("H144").Value = 3

$ws.Range("A145").Value = "Togo"
$ws.Range("B145").Value = 36
$ws.Range("D145").Value = 10
$ws.Range("E145").Value = 24
$ws.Range("H145").Value = 2

$ws.Range("A146").Value = "Guam"
$ws.Range("D146").Value = 0
$ws.Range("E146").Value = 31
$ws.Range("H146").Value = 1

$ws.Range("A147").Value = "Bermudas"
$ws.Range("B147").Value = 32
$ws.Range("D147").Value = 10
$ws.Range("E147").Value = 22
$ws.Range("H147").Value = 0

$ws.Range("A165").Value = "Mongolia"

$ws.Range("A166").Value = "Namibia"

$ws.Range("A171").Value = "Mozambique"

$ws.Range("A172").Value = "Seychelles"

$ws.Range("A173").Value = "Libia"

$ws.Range("A174").Value = "Granada"
$ws.Range("C174").Value = 1

$ws.Range("A175").Value = "Laos"
$ws.Range("C175").Value = 0

$ws.Range("A176").Value = "Surinam"

$ws.Range("A185").Value = "Antigua y Barbuda"

$ws.Range("A186").Value = "Fiyi"
$ws.Range("C186").Value = 2

$ws.Range("A187").Value = "Republica del Chad"
$ws.Range("C187").Value = 0

$ws.Range("A196").Value = "Somalia"
$ws.Range("D196").Value = 1
$ws.Range("H196").Value = 0

$ws.Range("A197").Value = "Nicaragua"
$ws.Range("D197").Value = 0
$ws.Range("H197").Value = 1

$ws.Range("A200").Value = "Malaui"
$ws.Range("C200").Value = 3

$ws.Range("A201").Value = "Republica de Africa Central"
$ws.Range("C201").Value = 0

$ws.Range("A202").Value = "Belice"

$ws.Range("A203").Value = "Islas Virgenes Britanicas"
$ws.Range("C203").Value = 0

$ws.Range("A204").Value = "Burundi"
$ws.Range("C204").Value = 1

$ws.Range("A205").Value = "Anguila"
$ws.Range("B205").Value = 3
$ws.Range("C205").Value = 1
$ws.Range("E205").Value = 3

$ws.Range("A207").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("D207").Value = 0
$ws.Range("E207").Value = 2

$ws.Range("A208").Value = "San Vicente y las Granadinas"
$ws.Range("B208").Value = 2
$ws.Range("D208").Value = 1

# Add new row 210 (Papua Nueva Guinea entry, previously at row 208 position in data)
$ws.Range("A210").Value = "Papua Nueva Guinea"
$ws.Range("B210").Value = 1
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 0
$ws.Range("E210").Value = 1
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0
